$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in columns A, B, E, F, G, H between row 16 and row 18
$cols = @("A", "B", "E", "F", "G", "H")

foreach ($col in $cols) {
    $cell16 = $ws.Range("$col" + "16")
    $cell18 = $ws.Range("$col" + "18")
    $v16 = $cell16.Value2
    $v18 = $cell18.Value2
    $cell16.Value2 = $v18
    $cell18.Value2 = $v16
}
